$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value forcing text storage so that numeric-looking
# strings (e.g. "61.16") are not silently converted/rounded into floating
# point numbers by Excel's type inference.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.724.55"
$ws.Range("E2").Value = "  -0.88%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.544.52"
$ws.Range("E3").Value = "  -1.38%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.29%  "

# Row 5 - BNB
Set-TextValue "D5" "206.12"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.21%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  -3.21%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.20%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0581"
$ws.Range("E10").Value = "  -0.72%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.71%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.763.97"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.547.68"
$ws.Range("E13").Value = "  -1.19%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.16%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.88%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.711.68"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17 - Litecoin
Set-TextValue "D17" "61.16"
$ws.Range("E17").Value = "  -0.93%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "212.46"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.23"
$ws.Range("E20").Value = "  -1.66%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.26%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.98%  "

# Row 23 - Avalanche
Set-TextValue "D23" "8.96"
$ws.Range("E23").Value = "  -4.62%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.38%  "

# Row 25 - Monero
Set-TextValue "D25" "152.23"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "14.87"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - Cosmos
Set-TextValue "D27" "6.48"
$ws.Range("E27").Value = "  -2.99%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.28%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -0.98%  "

# Row 30 - was PancakeSwap, now Hedera
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0459"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31 - was Hedera, now PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.10"
$ws.Range("E31").Value = "  -1.55%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.32%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.347.57"
$ws.Range("E33").Value = "  -2.95%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "2.92"
$ws.Range("E34").Value = "  +0.35%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.92%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.94%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "0.932"
$ws.Range("E37").Value = "  -1.15%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.28%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.16%  "

# Row 40 - was FraxShare, now ARBITRUM
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.801"
$ws.Range("E40").Value = "  -1.00%  "

# Row 41 - was ARBITRUM, now FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.72"
$ws.Range("E41").Value = "  +5.48%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  -1.31%  "

# Row 43 - MXToken
Set-TextValue "D43" "2.18"
$ws.Range("E43").Value = "  -0.08%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  -3.37%  "

# Row 45 - Aave
Set-TextValue "D45" "62.59"
$ws.Range("E45").Value = "  -1.45%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.678.71"
$ws.Range("E46").Value = "  -1.35%  "

# Row 47 - mCoin
$ws.Range("E47").Value = "  -4.22%  "

# Row 48 - Quant
Set-TextValue "D48" "85.74"
$ws.Range("E48").Value = "  +0.46%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0507"
$ws.Range("E49").Value = "  +2.31%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₇0974"
$ws.Range("E50").Value = "  +0.25%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0952"
$ws.Range("E51").Value = "  +0.28%  "
